# Revert responsive design implementation
# Re-applies the pre-edit sensor rows: restores several trailing data rows
# that had been trimmed from each sheet and rolls back a batch of
# "responsive" timestamp adjustments on ROW02-FE-LIFTER.

$wb = $excel.ActiveWorkbook

$DATEFMT = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------
# Sheet "ROW35-FE-LIFTER": append rows 29-31 (clones of row 28's pattern
# with new timestamps).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")

$rows1 = @(
    @(45729.73239443287, "0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,", "0x01,0x90,", "0xd", 400, 568631262647113769549824.0, 400, 13),
    @(45729.73241640046, "0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,", "0x01,0x90,", "0xd", 400, 568631262647113769549824.0, 400, 13),
    @(45729.73243972223, "0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,", "0x01,0x90,", "0xd", 400, 568631262647113769549824.0, 400, 13)
)

$r = 29
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 1).NumberFormat = $DATEFMT
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $ws1.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "ROW35-MID-LIFTER": convert existing row 29's text timestamp into
# a numeric date, insert two new numeric rows (30-31), then re-append the
# original trailing text-timestamp row as row 32.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")

$ws2.Cells.Item(29, 1).Value = 45729.58037369213
$ws2.Cells.Item(29, 1).NumberFormat = $DATEFMT

$rows2 = @(
    @(45729.58039555555, "0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,", "0x01,0x86,", "0x4", 400, 568631262647113769549824.0, 390, 4),
    @(45729.58041870371, "0x01,0x90", "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,", "0x01,0x86,", "0x4", 400, 568631262647113769549824.0, 390, 4)
)

$r = 30
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 1).NumberFormat = $DATEFMT
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $ws2.Cells.Item($r, 8).Value = $row[7]
    $ws2.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

$ws2.Cells.Item(32, 1).Value = "2025-03-14 01:55:48"
$ws2.Cells.Item(32, 2).Value = "0x01,0x90"
$ws2.Cells.Item(32, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws2.Cells.Item(32, 4).Value = "0x01,0x86,"
$ws2.Cells.Item(32, 5).Value = "0x4"
$ws2.Cells.Item(32, 6).Value = 400
$ws2.Cells.Item(32, 7).Value = 568631262647113769549824.0
$ws2.Cells.Item(32, 8).Value = 390
$ws2.Cells.Item(32, 9).Value = 4

# ---------------------------------------------------------------------
# Sheet "ROW02-MID-LIFTER": append row 32 (text timestamp, clone of row
# 31's pattern).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ROW02-MID-LIFTER")

$ws3.Cells.Item(32, 1).Value = "2025-03-13 18:26:04"
$ws3.Cells.Item(32, 2).Value = "0x01,0x90"
$ws3.Cells.Item(32, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws3.Cells.Item(32, 4).Value = "0x01,0x90,"
$ws3.Cells.Item(32, 5).Value = "0x3"
$ws3.Cells.Item(32, 6).Value = 400
$ws3.Cells.Item(32, 7).Value = 568631262647113769549824.0
$ws3.Cells.Item(32, 8).Value = 400
$ws3.Cells.Item(32, 9).Value = 3

# ---------------------------------------------------------------------
# Sheet "ROW02-FE-LIFTER": roll back the 59-second timestamp adjustment
# that had been applied to rows 2-31's column A, then append new row 32.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("ROW02-FE-LIFTER")

$newDates4 = @(
    45725.28336619213,
    45725.28338791666,
    45725.28341164352,
    45725.78350880787,
    45725.78353120371,
    45725.78355416666,
    45726.28365108796,
    45726.28367324074,
    45726.28369701389,
    45726.78379283565,
    45726.78381515046,
    45726.78383841435,
    45727.28393501158,
    45727.28395739583,
    45727.28398025463,
    45727.78407800926,
    45727.78409938658,
    45727.78412253472,
    45728.28421967592,
    45728.28424142361,
    45728.28426457176,
    45728.32193778935,
    45728.32196105324,
    45728.32198430556,
    45728.82208192129,
    45728.82210399306,
    45728.8221271412,
    45729.32222401621,
    45729.32224596065,
    45729.32226951389
)

$r = 2
foreach ($d in $newDates4) {
    $ws4.Cells.Item($r, 1).Value = $d
    $ws4.Cells.Item($r, 1).NumberFormat = $DATEFMT
    $r = $r + 1
}

$ws4.Cells.Item(32, 1).Value = 45729.82226851852
$ws4.Cells.Item(32, 1).NumberFormat = $DATEFMT
$ws4.Cells.Item(32, 2).Value = "0x01,0x90"
$ws4.Cells.Item(32, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws4.Cells.Item(32, 4).Value = "0x01,0x90,"
$ws4.Cells.Item(32, 5).Value = "0x3"
$ws4.Cells.Item(32, 6).Value = 400
$ws4.Cells.Item(32, 7).Value = 985046333984776143241216.0
$ws4.Cells.Item(32, 8).Value = 400
$ws4.Cells.Item(32, 9).Value = 3
